# Applies the commit: insert one new weekly record row into the daily
# price series for "Pepino ensalada" (Terminal Hortofrutícola Agro Chillán),
# pushing the existing rows 157..245 down to 158..246 and adding the new
# observation at row 157.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Sheet1")

# Insert a new row before row 157; this shifts rows 157-245 down to 158-246
# and the sheet's used range grows from A1:R245 to A1:R246.
$ws.Rows.Item(157).Insert()

# Populate the newly inserted row with the new weekly observation.
$newRow = 157
$ws.Cells.Item($newRow, 1).Value = 7
$ws.Cells.Item($newRow, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item($newRow, 3).Value = "Ñuble"
$ws.Cells.Item($newRow, 4).Value = 44777
$ws.Cells.Item($newRow, 5).Value = 16
$ws.Cells.Item($newRow, 6).Value = 100112043
$ws.Cells.Item($newRow, 7).Value = "Pepino ensalada"
$ws.Cells.Item($newRow, 8).Value = "Sin especificar"
$ws.Cells.Item($newRow, 9).Value = "Primera"
$ws.Cells.Item($newRow, 10).Value = 100
$ws.Cells.Item($newRow, 11).Value = 19000
$ws.Cells.Item($newRow, 12).Value = 20000
$ws.Cells.Item($newRow, 13).Value = 19500
$ws.Cells.Item($newRow, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item($newRow, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item($newRow, 16).Value = 325
$ws.Cells.Item($newRow, 17).Value = 60
$ws.Cells.Item($newRow, 18).Value = "Hortaliza"
